# Natmi following Dr Hou advice
# Recompute Icam1-Itgax ligand-receptor edges with min_expressing set to 3
# (was 1), which produces 8 sending/target cluster pairs (ECs/FAPs/M2/sCs -> ECs/M2)
# instead of the previous 4 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Icam1"
$ws.Cells.Item(2, 3).Value = "Itgax"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 144.843106
$ws.Cells.Item(2, 8).Value = 434.529318
$ws.Cells.Item(2, 9).Value = 0.6517202749316883
$ws.Cells.Item(2, 10).Value = 0.6517202749316884
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 15.04425566666667
$ws.Cells.Item(2, 14).Value = 45.132767
$ws.Cells.Item(2, 15).Value = 0.3565971854932169
$ws.Cells.Item(2, 16).Value = 0.356597185493217
$ws.Cells.Item(2, 17).Value = 2179.056718218101
$ws.Cells.Item(2, 18).Value = 19611.51046396291
$ws.Cells.Item(2, 19).Value = 0.2324016157695055
$ws.Cells.Item(2, 20).Value = 0.2324016157695056

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Icam1"
$ws.Cells.Item(3, 3).Value = "Itgax"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 144.843106
$ws.Cells.Item(3, 8).Value = 434.529318
$ws.Cells.Item(3, 9).Value = 0.6517202749316883
$ws.Cells.Item(3, 10).Value = 0.6517202749316884
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 27.144119
$ws.Cells.Item(3, 14).Value = 81.432357
$ws.Cells.Item(3, 15).Value = 0.643402814506783
$ws.Cells.Item(3, 16).Value = 0.643402814506783
$ws.Cells.Item(3, 17).Value = 3931.638505593614
$ws.Cells.Item(3, 18).Value = 35384.74655034253
$ws.Cells.Item(3, 19).Value = 0.4193186591621827
$ws.Cells.Item(3, 20).Value = 0.4193186591621827

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Icam1"
$ws.Cells.Item(4, 3).Value = "Itgax"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 33.793597
$ws.Cells.Item(4, 8).Value = 101.380791
$ws.Cells.Item(4, 9).Value = 0.1520539909422453
$ws.Cells.Item(4, 10).Value = 0.1520539909422453
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 15.04425566666667
$ws.Cells.Item(4, 14).Value = 45.132767
$ws.Cells.Item(4, 15).Value = 0.3565971854932169
$ws.Cells.Item(4, 16).Value = 0.356597185493217
$ws.Cells.Item(4, 17).Value = 508.3995131642997
$ws.Cells.Item(4, 18).Value = 4575.595618478697
$ws.Cells.Item(4, 19).Value = 0.05422202521301577
$ws.Cells.Item(4, 20).Value = 0.05422202521301578

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Icam1"
$ws.Cells.Item(5, 3).Value = "Itgax"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 33.793597
$ws.Cells.Item(5, 8).Value = 101.380791
$ws.Cells.Item(5, 9).Value = 0.1520539909422453
$ws.Cells.Item(5, 10).Value = 0.1520539909422453
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 27.144119
$ws.Cells.Item(5, 14).Value = 81.432357
$ws.Cells.Item(5, 15).Value = 0.643402814506783
$ws.Cells.Item(5, 16).Value = 0.643402814506783
$ws.Cells.Item(5, 17).Value = 917.297418406043
$ws.Cells.Item(5, 18).Value = 8255.676765654387
$ws.Cells.Item(5, 19).Value = 0.09783196572922953
$ws.Cells.Item(5, 20).Value = 0.09783196572922953

$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Icam1"
$ws.Cells.Item(6, 3).Value = "Itgax"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 41.33760066666667
$ws.Cells.Item(6, 8).Value = 124.012802
$ws.Cells.Item(6, 9).Value = 0.1859981687460937
$ws.Cells.Item(6, 10).Value = 0.1859981687460937
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 15.04425566666667
$ws.Cells.Item(6, 14).Value = 45.132767
$ws.Cells.Item(6, 15).Value = 0.3565971854932169
$ws.Cells.Item(6, 16).Value = 0.356597185493217
$ws.Cells.Item(6, 17).Value = 621.8934330759038
$ws.Cells.Item(6, 18).Value = 5597.040897683134
$ws.Cells.Item(6, 19).Value = 0.06632642348174944
$ws.Cells.Item(6, 20).Value = 0.06632642348174946

$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Icam1"
$ws.Cells.Item(7, 3).Value = "Itgax"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 41.33760066666667
$ws.Cells.Item(7, 8).Value = 124.012802
$ws.Cells.Item(7, 9).Value = 0.1859981687460937
$ws.Cells.Item(7, 10).Value = 0.1859981687460937
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 27.144119
$ws.Cells.Item(7, 14).Value = 81.432357
$ws.Cells.Item(7, 15).Value = 0.643402814506783
$ws.Cells.Item(7, 16).Value = 0.643402814506783
$ws.Cells.Item(7, 17).Value = 1122.072751670479
$ws.Cells.Item(7, 18).Value = 10098.65476503431
$ws.Cells.Item(7, 19).Value = 0.1196717452643443
$ws.Cells.Item(7, 20).Value = 0.1196717452643443

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Icam1"
$ws.Cells.Item(8, 3).Value = "Itgax"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.273049333333333
$ws.Cells.Item(8, 8).Value = 6.819148
$ws.Cells.Item(8, 9).Value = 0.01022756537997252
$ws.Cells.Item(8, 10).Value = 0.01022756537997253
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 15.04425566666667
$ws.Cells.Item(8, 14).Value = 45.132767
$ws.Cells.Item(8, 15).Value = 0.3565971854932169
$ws.Cells.Item(8, 16).Value = 0.356597185493217
$ws.Cells.Item(8, 17).Value = 34.19633531361288
$ws.Cells.Item(8, 18).Value = 307.767017822516
$ws.Cells.Item(8, 19).Value = 0.003647121028946066
$ws.Cells.Item(8, 20).Value = 0.003647121028946067

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Icam1"
$ws.Cells.Item(9, 3).Value = "Itgax"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.273049333333333
$ws.Cells.Item(9, 8).Value = 6.819148
$ws.Cells.Item(9, 9).Value = 0.01022756537997252
$ws.Cells.Item(9, 10).Value = 0.01022756537997253
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 27.144119
$ws.Cells.Item(9, 14).Value = 81.432357
$ws.Cells.Item(9, 15).Value = 0.643402814506783
$ws.Cells.Item(9, 16).Value = 0.643402814506783
$ws.Cells.Item(9, 17).Value = 61.69992159687067
$ws.Cells.Item(9, 18).Value = 555.299294371836
$ws.Cells.Item(9, 19).Value = 0.006580444351026458
$ws.Cells.Item(9, 20).Value = 0.00658044435102646

